# Commit: "Tue, Jun 30, 2020  4:05:45 PM"
#
# 1) Three tables switch their table style from the built-in
#    "No Style, Table Grid" id to the built-in "No Style, No Grid" id.
# 2) The presentation's applied theme (design) colour palette changes
#    from the custom "Red Violet" / Integral scheme to the stock
#    Office colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table that currently uses the old table style ---
$oldTableStyle = "{FDF85083-7E51-4256-89D2-16EECFA9340C}"
$newTableStyle = "{1925C07D-70DC-4C7A-AEC3-A6D984D03569}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldTableStyle) {
                $tbl.ApplyStyle($newTableStyle)
            }
        }
    }
}

# --- 2. Swap the theme colour scheme to the stock "Office" palette ---
# Helper: pack R,G,B into the BGR-ordered long that PowerPoint's
# ColorFormat/ThemeColor .RGB property uses.
function ToOleColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Index order exposed by ThemeColorScheme: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - matching the <a:clrScheme> child order.
$officeColors = @(
    @(0, 0, 0),        # 1  dk1      000000
    @(255, 255, 255),  # 2  lt1      FFFFFF
    @(68, 84, 106),    # 3  dk2      44546A
    @(231, 230, 230),  # 4  lt2      E7E6E6
    @(91, 155, 213),   # 5  accent1  5B9BD5
    @(237, 125, 49),   # 6  accent2  ED7D31
    @(165, 165, 165),  # 7  accent3  A5A5A5
    @(255, 192, 0),    # 8  accent4  FFC000
    @(68, 114, 196),   # 9  accent5  4472C4
    @(112, 173, 71),   # 10 accent6  70AD47
    @(5, 99, 193),     # 11 hlink    0563C1
    @(149, 79, 114)    # 12 folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $themeColors.Item($i).RGB = ToOleColor $rgb[0] $rgb[1] $rgb[2]
}
